# i18n.xlsx ("removed excel file i18n"):
# Six msgid/en translation-table rows for the BBMRI-Negotiator "directory
# export" feature were dropped from the i18nstrings sheet:
#   dataexplorer_directory_export_button           / Go to sample / data negotiation
#   dataexplorer_directory_export_dialog_yes        / Yes, Send to Negotiator
#   dataexplorer_directory_export_dialog_no         / No, I want to keep filtering
#   dataexplorer_directory_export_dialog_title      / Send request to the BBMRI Negotiator?
#   dataexplorer_directory_export_dialog_message    / Your current selection of biobanks ...
#   dataexplorer_directory_export_no_filters        / Please filter the collections ...
# In the original file these occupy rows 13:18 (A=msgid col, B=en col).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole rows: everything below shifts up (old row 19 -> new row
# 13, ..., old row 43 -> new row 37), the dimension shrinks from B43 to
# B37, and the shared-strings table is pruned to only the strings still
# referenced (count 86->74, uniqueCount 85->73) -- all automatically.
$ws.Range("A13:B18").EntireRow.Delete()

# Match the editor's final cursor/scroll position.
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("C25:C26").Select()
